$d = $word.ActiveDocument

$replacements = @(
    @("78×60=4680", "55×97=5335"),
    @("30×84=2520", "86×34=2924"),
    @("84×55=4620", "48×31=1488"),
    @("96×39=3744", "36×34=1224"),
    @("38×39=1482", "16×23=368"),
    @("86×95=8170", "60×69=4140"),
    @("28×30=840",  "12×32=384"),
    @("67×29=1943", "63×51=3213"),
    @("91×32=2912", "65×36=2340"),
    @("72×72=5184", "75×49=3675"),
    @("67×13=871",  "38×42=1596"),
    @("66×94=6204", "80×48=3840"),
    @("47×84=3948", "77×19=1463"),
    @("31×22=682",  "14×30=420"),
    @("44×20=880",  "38×11=418"),
    @("81×97=7857", "11×42=462"),
    @("98×36=3528", "18×66=1188"),
    @("81×88=7128", "27×30=810"),
    @("59×36=2124", "20×26=520"),
    @("75×87=6525", "94×66=6204"),
    @("78×46=3588", "50×37=1850"),
    @("24×75=1800", "22×93=2046"),
    @("57×30=1710", "15×66=990"),
    @("16×27=432",  "91×45=4095"),
    @("51×95=4845", "22×29=638")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
